# Correction to makeham_law_mortality.py: update premium1 (column E) values
# for rows 2-22 on the IAx_lifeTables3 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IAx_lifeTables3")

$newValues = @{
    2  = 2906.415819431462
    3  = 3503.120089918613
    4  = 4212.486766063165
    5  = 5019.805981111715
    6  = 5917.479452045765
    7  = 6881.165390848602
    8  = 7741.474704106309
    9  = 2369.727268068475
    10 = 2858.180025836033
    11 = 3433.328240030961
    12 = 4109.572641552428
    13 = 4897.597309379925
    14 = 5776.868550542882
    15 = 6673.746067758563
    16 = 2870.401035842662
    17 = 3415.861063987529
    18 = 4080.743127840228
    19 = 4843.955067316742
    20 = 5659.079067217303
    21 = 6482.947163918746
    22 = 7212.195501635304
}

foreach ($row in $newValues.Keys) {
    $ws.Range("E$row").Value = $newValues[$row]
}
